# Scheduled market-data refresh: updates currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) for a set of Leve rows across all eight job
# sheets, reflecting freshly pulled Universalis price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 18518754
$ws.Range("J6").Value = 349.1111
$ws.Range("L6").Value = 1047.3333
$ws.Range("N6").Value = -1271.3333

$ws.Range("H125").Value = 2935.1
$ws.Range("I125").Value = 784.3333
$ws.Range("J125").Value = 3856.8572
$ws.Range("K125").Value = 7058.9997
$ws.Range("L125").Value = 34711.7148
$ws.Range("M125").Value = -4598.9997
$ws.Range("N125").Value = -39631.7148

$ws.Range("H137").Value = 435518.84
$ws.Range("I137").Value = 1031837.3
$ws.Range("K137").Value = 3095511.9
$ws.Range("M137").Value = -3092961.9

$ws.Range("H138").Value = 10354.939
$ws.Range("J138").Value = 9788.52
$ws.Range("L138").Value = 29365.56
$ws.Range("N138").Value = -39645.56

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 39428.758
$ws.Range("I2").Value = 5747.524
$ws.Range("K2").Value = 5747.524
$ws.Range("M2").Value = -5634.524

$ws.Range("H61").Value = 2777.7778
$ws.Range("I61").Value = 818.4545000000001
$ws.Range("J61").Value = 5856.7144
$ws.Range("K61").Value = 818.4545000000001
$ws.Range("L61").Value = 5856.7144
$ws.Range("M61").Value = -606.4545000000001
$ws.Range("N61").Value = -6280.7144

$ws.Range("H110").Value = 1052
$ws.Range("I110").Value = 461.84616
$ws.Range("K110").Value = 461.84616
$ws.Range("M110").Value = 1583.15384

$ws.Range("H116").Value = 39428.758
$ws.Range("I116").Value = 5747.524
$ws.Range("K116").Value = 5747.524
$ws.Range("M116").Value = -3453.524

$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("M117").ClearContents()

$ws.Range("H136").Value = 2777.7778
$ws.Range("I136").Value = 818.4545000000001
$ws.Range("J136").Value = 5856.7144
$ws.Range("K136").Value = 2455.3635
$ws.Range("L136").Value = 17570.1432
$ws.Range("M136").Value = 94.63649999999961
$ws.Range("N136").Value = -22670.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 39428.758
$ws.Range("I3").Value = 5747.524
$ws.Range("K3").Value = 5747.524
$ws.Range("M3").Value = -5633.524

$ws.Range("H20").Value = 2621.3914
$ws.Range("I20").Value = 1506.1538
$ws.Range("K20").Value = 1506.1538
$ws.Range("M20").Value = -1259.1538

$ws.Range("H54").Value = 17737.143
$ws.Range("I54").Value = 14832.2
$ws.Range("K54").Value = 14832.2
$ws.Range("M54").Value = -14348.2

$ws.Range("H132").Value = 89999
$ws.Range("J132").Value = 89999
$ws.Range("L132").Value = 89999
$ws.Range("N132").Value = -100119

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1931.4615
$ws.Range("I16").Value = 2416.2856
$ws.Range("K16").Value = 2416.2856
$ws.Range("M16").Value = -2129.2856

$ws.Range("H59").Value = 56111.89
$ws.Range("J59").Value = 57625.875
$ws.Range("L59").Value = 57625.875
$ws.Range("N59").Value = -59915.875

$ws.Range("H113").Value = 1931.4615
$ws.Range("I113").Value = 2416.2856
$ws.Range("K113").Value = 2416.2856
$ws.Range("M113").Value = -246.2856000000002

$ws.Range("H132").Value = 19602.035
$ws.Range("I132").Value = 6248.5386
$ws.Range("K132").Value = 18745.6158
$ws.Range("M132").Value = -16215.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 3699.75
$ws.Range("I116").Value = 1599.6666
$ws.Range("J116").Value = 10000
$ws.Range("K116").Value = 4798.9998
$ws.Range("L116").Value = 30000
$ws.Range("M116").Value = -1356.9998
$ws.Range("N116").Value = -36884

$ws.Range("H123").Value = 4032.5
$ws.Range("I123").Value = 2265
$ws.Range("K123").Value = 6795
$ws.Range("M123").Value = -4345

$ws.Range("H132").Value = 8799599
$ws.Range("I132").Value = 1092.7142
$ws.Range("J132").Value = 13932062
$ws.Range("K132").Value = 9834.427799999999
$ws.Range("L132").Value = 125388558
$ws.Range("M132").Value = -7304.427799999999
$ws.Range("N132").Value = -125393618

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 106.35714
$ws.Range("I2").Value = 122.09091
$ws.Range("J2").Value = 48.666668
$ws.Range("K2").Value = 122.09091
$ws.Range("L2").Value = 48.666668
$ws.Range("M2").Value = -9.090909999999994
$ws.Range("N2").Value = -274.666668

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H70").Value = 6120.8887
$ws.Range("I70").Value = 6560.2
$ws.Range("K70").Value = 6560.2
$ws.Range("M70").Value = -6290.2

$ws.Range("H73").Value = 6120.8887
$ws.Range("I73").Value = 6560.2
$ws.Range("K73").Value = 6560.2
$ws.Range("M73").Value = -5624.2

$ws.Range("H102").Value = 16836.684
$ws.Range("I102").Value = 17817.53
$ws.Range("K102").Value = 17817.53
$ws.Range("M102").Value = -16195.53

$ws.Range("H107").Value = 711
$ws.Range("J107").Value = 717.5714
$ws.Range("L107").Value = 717.5714
$ws.Range("N107").Value = -4557.5714

$ws.Range("H110").Value = 140000
$ws.Range("J110").Value = 140000
$ws.Range("L110").Value = 140000
$ws.Range("N110").Value = -148180

$ws.Range("H126").Value = 13409.538
$ws.Range("I126").Value = 18999.9
$ws.Range("J126").Value = 9915.5625
$ws.Range("K126").Value = 56999.7
$ws.Range("L126").Value = 29746.6875
$ws.Range("M126").Value = -54529.7
$ws.Range("N126").Value = -34686.6875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 22505.5
$ws.Range("I29").Value = 22505.5
$ws.Range("K29").Value = 22505.5
$ws.Range("M29").Value = -22210.5

$ws.Range("H61").Value = 1435.1428
$ws.Range("I61").Value = 1191
$ws.Range("K61").Value = 1191
$ws.Range("M61").Value = -989

$ws.Range("H82").Value = 1749.3055
$ws.Range("I82").Value = 2160
$ws.Range("K82").Value = 2160
$ws.Range("M82").Value = -1799

$ws.Range("H85").Value = 1749.3055
$ws.Range("I85").Value = 2160
$ws.Range("K85").Value = 2160
$ws.Range("M85").Value = -912

$ws.Range("H100").Value = 6728.143
$ws.Range("I100").Value = 7526.727
$ws.Range("K100").Value = 7526.727
$ws.Range("M100").Value = -6985.727

$ws.Range("H113").Value = 1435.1428
$ws.Range("I113").Value = 1191
$ws.Range("K113").Value = 1191
$ws.Range("M113").Value = 979

$ws.Range("H132").Value = 17127.54
$ws.Range("I132").Value = 20965.8
$ws.Range("J132").Value = 4333.3335
$ws.Range("K132").Value = 62897.39999999999
$ws.Range("L132").Value = 13000.0005
$ws.Range("M132").Value = -60367.39999999999
$ws.Range("N132").Value = -18060.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 15013.75
$ws.Range("I32").Value = 18013
$ws.Range("J32").Value = 12014.5
$ws.Range("K32").Value = 18013
$ws.Range("L32").Value = 12014.5
$ws.Range("M32").Value = -17696
$ws.Range("N32").Value = -12648.5

$ws.Range("H41").Value = 11957
$ws.Range("J41").Value = 11199.429
$ws.Range("L41").Value = 11199.429
$ws.Range("N41").Value = -11979.429

$ws.Range("H62").Value = 107628.44
$ws.Range("I62").Value = 186881.31
$ws.Range("K62").Value = 186881.31
$ws.Range("M62").Value = -186257.31

$ws.Range("H65").Value = 107628.44
$ws.Range("I65").Value = 186881.31
$ws.Range("K65").Value = 934406.55
$ws.Range("M65").Value = -931286.55

$ws.Range("H122").Value = 4591.7
$ws.Range("I122").Value = 2950.7778
$ws.Range("K122").Value = 8852.3334
$ws.Range("M122").Value = -6402.3334

$ws.Range("H125").Value = 61249.75
$ws.Range("J125").Value = 61249.75
$ws.Range("L125").Value = 61249.75
$ws.Range("N125").Value = -71089.75

$ws.Range("H126").Value = 21856.166
$ws.Range("I126").Value = 40360.363
$ws.Range("J126").Value = 6198.769
$ws.Range("K126").Value = 121081.089
$ws.Range("L126").Value = 18596.307
$ws.Range("M126").Value = -118611.089
$ws.Range("N126").Value = -23536.307
